$wb = $excel.ActiveWorkbook

# ALC!row2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2478.1
$ws.Range("I2").Value = 1572.375
$ws.Range("K2").Value = 1572.375
$ws.Range("M2").Value = -1459.375

# ALC!row6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 181.66667
$ws.Range("I6").Value = 45
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 135
$ws.Range("L6").Value = 750
$ws.Range("M6").Value = -23
$ws.Range("N6").Value = -974

# ALC!row39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 202.5
$ws.Range("I39").Value = 69.22221999999999
$ws.Range("J39").Value = 602.3333
$ws.Range("K39").Value = 207.66666
$ws.Range("L39").Value = 1806.9999
$ws.Range("M39").Value = 88.33334000000002
$ws.Range("N39").Value = -2398.9999

# ALC!row98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1244.3125
$ws.Range("I98").Value = 1175.75
$ws.Range("J98").Value = 1450
$ws.Range("K98").Value = 1175.75
$ws.Range("L98").Value = 1450
$ws.Range("M98").Value = 322.25
$ws.Range("N98").Value = -4446

# ALC!row113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2694.5
$ws.Range("I113").Value = 2074.8333
$ws.Range("J113").Value = 3066.3
$ws.Range("K113").Value = 2074.8333
$ws.Range("L113").Value = 3066.3
$ws.Range("M113").Value = 1179.1667
$ws.Range("N113").Value = -9574.299999999999

# ALC!row115
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 1770.238
$ws.Range("I115").Value = 1358.3334
$ws.Range("J115").Value = 2800
$ws.Range("K115").Value = 4075.0002
$ws.Range("L115").Value = 8400
$ws.Range("M115").Value = -2508.0002
$ws.Range("N115").Value = -11534

# ALC!row122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1244.3125
$ws.Range("I122").Value = 1175.75
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 3527.25
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -1077.25
$ws.Range("N122").Value = -9250

# ARM!row2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1736.0834
$ws.Range("J2").Value = 1066.6666
$ws.Range("L2").Value = 1066.6666
$ws.Range("N2").Value = -1292.6666

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4022.225
$ws.Range("I74").Value = 1701.9354
$ws.Range("J74").Value = 12014.333
$ws.Range("K74").Value = 1701.9354
$ws.Range("L74").Value = 12014.333
$ws.Range("M74").Value = -827.9354000000001
$ws.Range("N74").Value = -13762.333

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4022.225
$ws.Range("I77").Value = 1701.9354
$ws.Range("J77").Value = 12014.333
$ws.Range("K77").Value = 8509.677
$ws.Range("L77").Value = 60071.665
$ws.Range("M77").Value = -4141.677
$ws.Range("N77").Value = -68807.66500000001

# ARM!row88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 6500
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -3812

# ARM!row91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 6500
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -5808

# ARM!row116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1736.0834
$ws.Range("J116").Value = 1066.6666
$ws.Range("L116").Value = 1066.6666
$ws.Range("N116").Value = -5654.6666

# BSM!row3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1736.0834
$ws.Range("J3").Value = 1066.6666
$ws.Range("L3").Value = 1066.6666
$ws.Range("N3").Value = -1294.6666

# BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1616.6666
$ws.Range("I86").Value = 1615.1515
$ws.Range("K86").Value = 1615.1515
$ws.Range("M86").Value = -492.1514999999999

# BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1616.6666
$ws.Range("I89").Value = 1615.1515
$ws.Range("K89").Value = 8075.7575
$ws.Range("M89").Value = -2459.7575

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2946.7046
$ws.Range("I134").Value = 2021.35
$ws.Range("J134").Value = 3717.8333
$ws.Range("K134").Value = 6064.049999999999
$ws.Range("L134").Value = 11153.4999
$ws.Range("M134").Value = -3529.049999999999
$ws.Range("N134").Value = -16223.4999

# CUL!row44
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 818.2
$ws.Range("J44").Value = 830.3333
$ws.Range("L44").Value = 2490.9999
$ws.Range("N44").Value = -3286.9999

# CUL!row118
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 2190.4211
$ws.Range("I118").Value = 506
$ws.Range("J118").Value = 3706.4
$ws.Range("K118").Value = 1518
$ws.Range("L118").Value = 11119.2
$ws.Range("M118").Value = -275
$ws.Range("N118").Value = -13605.2

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 20331.64
$ws.Range("I131").Value = 2790
$ws.Range("J131").Value = 22723.682
$ws.Range("K131").Value = 8370
$ws.Range("L131").Value = 68171.046
$ws.Range("M131").Value = -3330
$ws.Range("N131").Value = -78251.046

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6269.0586
$ws.Range("I70").Value = 5877.091
$ws.Range("J70").Value = 6456.522
$ws.Range("K70").Value = 5877.091
$ws.Range("L70").Value = 6456.522
$ws.Range("M70").Value = -5607.091
$ws.Range("N70").Value = -6996.522

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6269.0586
$ws.Range("I73").Value = 5877.091
$ws.Range("J73").Value = 6456.522
$ws.Range("K73").Value = 5877.091
$ws.Range("L73").Value = 6456.522
$ws.Range("M73").Value = -4941.091
$ws.Range("N73").Value = -8328.522000000001

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3018.8125
$ws.Range("I80").Value = 2335
$ws.Range("J80").Value = 3176.6155
$ws.Range("K80").Value = 2335
$ws.Range("L80").Value = 3176.6155
$ws.Range("M80").Value = -1337
$ws.Range("N80").Value = -5172.6155

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3018.8125
$ws.Range("I83").Value = 2335
$ws.Range("J83").Value = 3176.6155
$ws.Range("K83").Value = 11675
$ws.Range("L83").Value = 15883.0775
$ws.Range("M83").Value = -6683
$ws.Range("N83").Value = -25867.0775

# LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6477.5938
$ws.Range("I122").Value = 6106.884
$ws.Range("J122").Value = 7236.6665
$ws.Range("K122").Value = 18320.652
$ws.Range("L122").Value = 21709.9995
$ws.Range("M122").Value = -15870.652
$ws.Range("N122").Value = -26609.9995

# WVR!row113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 8860.77
$ws.Range("I113").Value = 20598.2
$ws.Range("J113").Value = 1524.875
$ws.Range("K113").Value = 61794.60000000001
$ws.Range("L113").Value = 4574.625
$ws.Range("M113").Value = -59624.60000000001
$ws.Range("N113").Value = -8914.625
